# Horarios actualizados Linea 141 - 99
# Applies the scraped-data refresh to all three sheets of the workbook:
#  - "LP1912"      (sheet1): header text + full data rewrite (44 -> 58 rows)
#  - "LP1912-215"  (sheet2): header text + full data rewrite (4 -> 7 rows)
#  - "6203-6173"   (sheet3): header text only (data unchanged)

$wb = $excel.ActiveWorkbook

$sheet1Data = @(
    ,@("16:46:42","16:47","15_ABASTO",1,"LP1912")
    ,@("16:50:41","16:50","10_OLMOS",0,"LP1912")
    ,@("16:50:41","16:51","15_ABASTO",1,"LP1912")
    ,@("16:52:37","16:53","15_ABASTO",1,"LP1912")
    ,@("16:46:42","16:53","10_OLMOS",7,"LP1912")
    ,@("16:46:42","16:56","215C_EL PATO",10,"LP1912")
    ,@("16:46:42","17:01","16_SANTA ANA",15,"LP1912")
    ,@("16:46:42","17:03","23_HERNANDEZ",17,"LP1912")
    ,@("16:46:42","17:04","14_ABASTO",18,"LP1912")
    ,@("16:46:42","17:07","15_ABASTO",21,"LP1912")
    ,@("16:46:42","17:13","23_HERNANDEZ",27,"LP1912")
    ,@("16:46:42","17:14","10_OLMOS",28,"LP1912")
    ,@("17:13:30","17:16","10_OLMOS",3,"LP1912")
    ,@("16:46:42","17:17","17_ROMERO",31,"LP1912")
    ,@("16:50:41","17:17","23_HERNANDEZ",27,"LP1912")
    ,@("16:52:37","17:20","23_HERNANDEZ",28,"LP1912")
    ,@("16:46:42","17:23","16_SANTA ANA",37,"LP1912")
    ,@("16:46:42","17:24","11_ETCHEVERRY",38,"LP1912")
    ,@("17:13:30","17:27","15_ABASTO",14,"LP1912")
    ,@("17:13:30","17:33","23_HERNANDEZ",20,"LP1912")
    ,@("16:50:41","17:34","10_OLMOS",44,"LP1912")
    ,@("16:46:42","17:35","16_P MOR-SANTA ANA",49,"LP1912")
    ,@("16:52:37","17:36","27_EL RETIRO",44,"LP1912")
    ,@("17:13:30","17:37","27_EL RETIRO",24,"LP1912")
    ,@("16:46:42","17:38","17X38_ROMERO",52,"LP1912")
    ,@("17:35:09","17:39","27_EL RETIRO",4,"LP1912")
    ,@("17:13:30","17:41","23_HERNANDEZ",28,"LP1912")
    ,@("16:46:42","17:44","215B_EL PATO",58,"LP1912")
    ,@("17:35:09","17:45","215B_EL PATO",10,"LP1912")
    ,@("16:50:41","17:47","16_SANTA ANA",57,"LP1912")
    ,@("16:46:42","17:48","27_EL RETIRO",62,"LP1912")
    ,@("16:50:41","17:49","27_EL RETIRO",59,"LP1912")
    ,@("16:46:42","17:50","215_EL PELIGRO",64,"LP1912")
    ,@("16:52:37","17:51","215_EL PELIGRO",59,"LP1912")
    ,@("17:35:09","18:00","16_SANTA ANA",25,"LP1912")
    ,@("16:46:42","18:02","17_ROMERO",76,"LP1912")
    ,@("17:35:09","18:03","23_HERNANDEZ",28,"LP1912")
    ,@("16:52:37","18:03","17_ROMERO",71,"LP1912")
    ,@("16:46:42","18:04","14_ABASTO",78,"LP1912")
    ,@("17:35:09","18:05","14_ABASTO",30,"LP1912")
    ,@("16:52:37","18:14","10_OLMOS",82,"LP1912")
    ,@("16:46:42","18:24","11_ETCHEVERRY",98,"LP1912")
    ,@("17:35:09","18:25","11_ETCHEVERRY",50,"LP1912")
    ,@("17:13:30","18:27","15_ABASTO",74,"LP1912")
    ,@("17:35:09","18:31","23_HERNANDEZ",56,"LP1912")
    ,@("16:46:42","18:34","14X44_ABASTO",108,"LP1912")
    ,@("16:46:42","18:38","17X38_ROMERO",112,"LP1912")
    ,@("17:13:30","18:41","14_ABASTO",88,"LP1912")
    ,@("16:46:42","18:41","16_P MOR-SANTA ANA",115,"LP1912")
    ,@("17:35:09","18:45","14_ABASTO",70,"LP1912")
    ,@("17:35:09","18:51","15_ABASTO",76,"LP1912")
    ,@("17:35:09","18:59","10_OLMOS",84,"LP1912")
    ,@("17:13:30","19:01","17_ROMERO",108,"LP1912")
    ,@("17:13:30","19:11","81_EL PELIGRO",118,"LP1912")
    ,@("17:35:09","19:19","27_EL RETIRO",104,"LP1912")
    ,@("17:35:09","19:21","215C_EL PATO",106,"LP1912")
    ,@("17:35:09","19:29","225_GOMEZ",114,"LP1912")
    ,@("17:35:09","19:31","215_EL PELIGRO",116,"LP1912")
)

$sheet2Data = @(
    ,@("16:46:42","16:56","215C_EL PATO",10,"LP1912")
    ,@("16:46:42","17:44","215B_EL PATO",58,"LP1912")
    ,@("17:35:09","17:45","215B_EL PATO",10,"LP1912")
    ,@("16:46:42","17:50","215_EL PELIGRO",64,"LP1912")
    ,@("16:52:37","17:51","215_EL PELIGRO",59,"LP1912")
    ,@("17:35:09","19:21","215C_EL PATO",106,"LP1912")
    ,@("17:35:09","19:31","215_EL PELIGRO",116,"LP1912")
)


$wsLP1912     = $wb.Worksheets.Item("LP1912")
$wsLP1912215  = $wb.Worksheets.Item("LP1912-215")
$ws6203       = $wb.Worksheets.Item("6203-6173")

# --- Sheet "LP1912" ---
$wsLP1912.Cells.Item(2, 1).Value = "Última actualización: 17:35:09"
$wsLP1912.Cells.Item(3, 1).Value = "Total filas: 58"

$r = 6
foreach ($row in $sheet1Data) {
    $wsLP1912.Cells.Item($r, 1).Value = $row[0]
    $wsLP1912.Cells.Item($r, 2).Value = $row[1]
    $wsLP1912.Cells.Item($r, 3).Value = $row[2]
    $wsLP1912.Cells.Item($r, 4).Value = $row[3]
    $wsLP1912.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# --- Sheet "LP1912-215" ---
$wsLP1912215.Cells.Item(2, 1).Value = "Última actualización: 17:35:09"
$wsLP1912215.Cells.Item(3, 1).Value = "Total filas: 7"

$r = 6
foreach ($row in $sheet2Data) {
    $wsLP1912215.Cells.Item($r, 1).Value = $row[0]
    $wsLP1912215.Cells.Item($r, 2).Value = $row[1]
    $wsLP1912215.Cells.Item($r, 3).Value = $row[2]
    $wsLP1912215.Cells.Item($r, 4).Value = $row[3]
    $wsLP1912215.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# --- Sheet "6203-6173" (header text only; data unchanged) ---
$ws6203.Cells.Item(2, 1).Value = "Última actualización: 17:35:09"

Write-Output "Updated LP1912 (58 rows), LP1912-215 (7 rows), and 6203-6173 header."
